# Set the "AB" values (column B) to 0 for years 1600-1900 (rows 2-302)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B302").Value = 0
